$wb = $excel.ActiveWorkbook

$oldName = "aabae1fa-5600-4ae5-b204-47f188ff0c50.md"
$newName = "ae013254-c540-4bff-a548-43c6ef4ab4af.md"

$oldStatus = "Ready for handoff"
$newStatus = "Handoff transform failed"

# --- Sheet "Overview" ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Hyperlinks.Item(1).TextToDisplay = $newName
$wsOverview.Range("A2").Value = $newName
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus

# --- Sheet "zh-cn" ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Hyperlinks.Item(1).TextToDisplay = $newName
$wsZh.Range("A2").Value = $newName
$wsZh.Range("B2").Value = $newStatus
$wsZh.Range("C2").ClearContents()
$wsZh.Range("D2").Value = "0001-01-01 00:00:00"
$wsZh.Range("H2").Value = "Ignored"

# --- Sheet "de-de" ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Hyperlinks.Item(1).TextToDisplay = $newName
$wsDe.Range("A2").Value = $newName
$wsDe.Range("B2").Value = $newStatus
$wsDe.Range("C2").ClearContents()
$wsDe.Range("D2").Value = "0001-01-01 00:00:00"
$wsDe.Range("H2").Value = "Ignored"
